$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.527.82"
$ws.Range("E2").Value = "  +3.23%  "
$ws.Range("D3").Value = "3.461.37"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.37%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.462.38"
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("E9").Value = "  +5.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("E11").Value = "  +6.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").Value = "4.061.92"
$ws.Range("E13").Value = "  +3.98%  "
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("E15").Value = "  +9.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").Value = "64.566.58"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "3.465.54"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "395.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.546"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +27.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.23%  "
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.99%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.71%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.44%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.74%  "
$ws.Range("E36").Value = "  +5.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0786"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.15%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "2.941.02"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.21%  "
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +25.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.863"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.15%  "
